$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-29 Saturday" "2024-06-30 Sunday"
Replace-Text "63×15=945" "96×61=5856"
Replace-Text "61×61=3721" "29×21=609"
Replace-Text "60×75=4500" "57×82=4674"
Replace-Text "23×93=2139" "38×18=684"
Replace-Text "62×88=5456" "48×48=2304"
Replace-Text "49×99=4851" "90×41=3690"
Replace-Text "26×12=312" "48×93=4464"
Replace-Text "99×73=7227" "97×73=7081"
Replace-Text "63×33=2079" "83×63=5229"
Replace-Text "34×75=2550" "72×93=6696"
Replace-Text "34×91=3094" "36×18=648"
Replace-Text "20×18=360" "26×98=2548"
Replace-Text "69×54=3726" "37×96=3552"
Replace-Text "91×77=7007" "19×46=874"
Replace-Text "37×38=1406" "42×90=3780"
Replace-Text "24×31=744" "25×71=1775"
Replace-Text "25×81=2025" "80×23=1840"
Replace-Text "54×72=3888" "32×88=2816"
Replace-Text "34×76=2584" "67×25=1675"
Replace-Text "66×55=3630" "36×56=2016"
Replace-Text "41×71=2911" "80×79=6320"
Replace-Text "66×11=726" "50×27=1350"
Replace-Text "35×35=1225" "74×85=6290"
Replace-Text "61×62=3782" "79×51=4029"
Replace-Text "36×13=468" "21×45=945"

Write-Output "Done"
